$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" (Strikeouts replaced by K count per commit message).
# Update values for rows 2-7 as per regenerated save_data.
$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 5
